$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(4264, 4676, 4769, 4769, 4858, 4858, 5129, 5129, 5129, 5129, 5316, 5316, 5316, 5316)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 3).Value = $values[$i]
}
